# Update the Rspo2-Lgr5 ligand-receptor pair table with refreshed TPM-based NATMI
# statistics: new sending/target cluster combinations involving "MuSCs" are
# added, the "Inflammatory-Mac" cluster is renamed to "Resolving-Mac", and all
# numeric edge-weight/specificity figures are recalculated for every row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column letters for the data columns (A = Sending cluster ... T = Edge total
# expression derived specificity), matching the header row already present
# in row 1.
$columns = @("A","B","C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T")

# Each inner array holds: row number, followed by the 20 column values in
# order (A..T) for that row.
$data = @(
    @(2, "FAPs","Rspo2","Lgr5","ECs",3,1,1.876175666666667,5.628527,0.9891011365778073,0.9927075980877177,2,0.6666666666666666,0.2158443333333333,0.647533,0.0304269596471442,0.0383225350174662,0.4049618859878889,3.644656973891,0.03009534036959741,0.03804307168982132),
    @(3, "FAPs","Rspo2","Lgr5","FAPs",3,1,1.876175666666667,5.628527,0.9891011365778073,0.9927075980877177,3,1,2.321162,6.963486,0.3272075825100088,0.4121155772426048,4.354907662791333,39.194168965122,0.3236413917575264,0.4091102648190395),
    @(4, "FAPs","Rspo2","Lgr5","MuSCs",3,1,1.876175666666667,5.628527,0.9891011365778073,0.9927075980877177,2,1,4.3846295,8.769259,0.6180887068188557,0.5189854958816471,8.226335175248833,49.358011051493,0.6113522424204373,0.515200845059033),
    @(5, "FAPs","Rspo2","Lgr5","Resolving-Mac",3,1,1.876175666666667,5.628527,0.9891011365778073,0.9927075980877177,2,0.6666666666666666,0.1722156666666667,0.516647,0.02427675102399122,0.03057639185828191,0.3231068432187778,2.907961588969,0.02401216203024616,0.03035341651982388),
    @(6, "MuSCs","Rspo2","Lgr5","ECs",1,0.5,0.0206735,0.041347,0.01089886342219268,0.007292401912282354,2,0.6666666666666666,0.2158443333333333,0.647533,0.0304269596471442,0.0383225350174662,0.004462257825166667,0.026773546951,0.0003316192775467926,0.000279463327644878),
    @(7, "MuSCs","Rspo2","Lgr5","FAPs",1,0.5,0.0206735,0.041347,0.01089886342219268,0.007292401912282354,3,1,2.321162,6.963486,0.3272075825100088,0.4121155772426048,0.047986542607,0.287919255642,0.003566190752482428,0.003005312423565317),
    @(8, "MuSCs","Rspo2","Lgr5","MuSCs",1,0.5,0.0206735,0.041347,0.01089886342219268,0.007292401912282354,2,1,4.3846295,8.769259,0.6180887068188557,0.5189854958816471,0.09064563796825001,0.362582551873,0.006736464398418403,0.003784650822614129),
    @(9, "MuSCs","Rspo2","Lgr5","Resolving-Mac",1,0.5,0.0206735,0.041347,0.01089886342219268,0.007292401912282354,2,0.6666666666666666,0.1722156666666667,0.516647,0.02427675102399122,0.03057639185828191,0.003560300584833333,0.021361803509,0.0002645889937450566,0.0002229753384580296)
)

foreach ($entry in $data) {
    $r = $entry[0]
    for ($i = 0; $i -lt $columns.Length; $i++) {
        $col = $columns[$i]
        $value = $entry[$i + 1]
        $ws.Range("$col$r").Value = $value
    }
}
